# Insert a new weekly record at row 250, pushing the existing rows
# 250-302 down to 251-303 (dimension grows from A1:R302 to A1:R303).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(250).Insert()

$ws.Cells.Item(250, 1).Value  = 7
$ws.Cells.Item(250, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(250, 3).Value  = "Ñuble"
$ws.Cells.Item(250, 4).Value  = 44995
$ws.Cells.Item(250, 5).Value  = 16
$ws.Cells.Item(250, 6).Value  = 100112043
$ws.Cells.Item(250, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(250, 8).Value  = "Sin especificar"
$ws.Cells.Item(250, 9).Value  = "Primera"
$ws.Cells.Item(250, 10).Value = 60
$ws.Cells.Item(250, 11).Value = 9000
$ws.Cells.Item(250, 12).Value = 9000
$ws.Cells.Item(250, 13).Value = 9000
$ws.Cells.Item(250, 14).Value = "$/caja 80 unidades"
$ws.Cells.Item(250, 15).Value = "Región del Maule"
$ws.Cells.Item(250, 16).Value = 112
$ws.Cells.Item(250, 17).Value = 80
$ws.Cells.Item(250, 18).Value = "Hortaliza"
